# This script applies the content edit described by the commit:
# the placeholder course-code strings "mmmm1".."mmmm5" (and the
# derived "موضوعات مختارة في mmmm") are renamed to "CSE1".."CSE5"
# ("موضوعات مختارة في CSE") throughout the "اسم المساق" (course name)
# column, plus the one "التخصص والمستوى" cell that spelled out the
# "Selected topics" course title.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "CSE2"
$ws.Range("D3").Value = "CSE2"
$ws.Range("D4").Value = "CSE1"
$ws.Range("D6").Value = "CSE1"
$ws.Range("D7").Value = "CSE2"
$ws.Range("D9").Value = "CSE5"
$ws.Range("D11").Value = "CSE4"
$ws.Range("D12").Value = "CSE3"
$ws.Range("D13").Value = "CSE5"
$ws.Range("C14").Value = "موضوعات مختارة في CSE"
$ws.Range("D14").Value = "CSE5"
$ws.Range("D16").Value = "CSE1"
$ws.Range("D28").Value = "CSE1"
$ws.Range("D29").Value = "CSE1"
$ws.Range("D40").Value = "CSE1"
$ws.Range("D42").Value = "CSE3"
$ws.Range("D43").Value = "CSE3"
$ws.Range("D47").Value = "CSE5"
$ws.Range("D48").Value = "CSE4"
$ws.Range("D52").Value = "CSE3"
$ws.Range("D55").Value = "CSE3"
$ws.Range("D56").Value = "CSE3"
$ws.Range("D70").Value = "CSE5"
$ws.Range("D71").Value = "CSE5"
$ws.Range("D72").Value = "CSE3"
$ws.Range("D74").Value = "CSE1"
$ws.Range("D75").Value = "CSE1"
$ws.Range("D78").Value = "CSE1"
$ws.Range("D80").Value = "CSE2"
$ws.Range("D101").Value = "CSE1"

# Move the visible/scrolled window and active-cell selection to match
# the author's final cursor position in the sheet.
$ws.Application.ActiveWindow.ScrollRow = 43
$ws.Range("D56").Select()
